# GaN 2022 Activity Guide (Hercules) - Japanese
#
# Replace the old "Perseids" campaign-period sentence (spread across many
# per-character/per-fragment runs, e.g. "201" + "8" + "年キャンペーン期間
# (対象：" + "ペルセウス" + ")：" + ... + "日") with the new, already
# translated Hercules campaign-period sentence, as a single plain run
# (matching the canonical OOXML, which carries no run-level formatting
# override for the replacement run).
#
# The fully concatenated old text reads:
#   2018年キャンペーン期間 (対象：ペルセウス)：、10月30日〜11月8日、11月29日〜12月8日
# and becomes:
#   年キャンペーン期間 対象：Hercules: 6月13〜22日、7月12〜21日、8月10〜19日
#
# This occurs 4 times in the document.

$d = $word.ActiveDocument

$oldText = "2018" + `
    "年キャンペーン期間 (対象：" + `
    "ペルセウス" + `
    ")：" + `
    "、" + `
    "10" + "月" + "30" + "日〜" + "11" + "月" + "8" + "日、" + `
    "11" + "月" + "29" + "日〜" + "12" + "月" + "8" + "日"

$newText = "年キャンペーン期間 対象：Hercules: 6月13〜22日、7月12〜21日、8月10〜19日"
$newLen = $newText.Length

$found = $true
$guard = 0
while ($found -and $guard -lt 20) {
    $guard = $guard + 1
    $range = $d.Content
    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
    if ($found) {
        $origStart = $range.Start
        $origEnd = $range.End

        # Insert the replacement as a brand-new, unformatted run right
        # before the matched (old) text. InsertBefore grows $range so it
        # covers "newText + oldText"; then drop just the old-text tail,
        # leaving a single clean <w:r><w:t>newText</w:t></w:r>.
        $range.InsertBefore($newText)

        $oldTail = $d.Range($origStart + $newLen, $origEnd + $newLen)
        $oldTail.Delete()
    }
}
